$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Norm_Sayisi (C) and Mevcut_Sayisi (D) values for several rows.
# Column E ("Fark") holds a shared formula (=D-C) and recalculates automatically.

# Row 5: Mühendis
$ws.Range("C5").Value = 23
$ws.Range("D5").Value = 22

# Row 6: Teknisyen/Tekniker
$ws.Range("C6").Value = 39

# Row 7: İşlem Yapan Personel
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 3

# Row 8: Arşiv Görevlisi
$ws.Range("C8").Value = 4

# Row 9: Destek Personeli
$ws.Range("C9").Value = 10

# Row 10: Şoför
$ws.Range("D10").Value = 0

$excel.Calculate()
$wb.Save()
